$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
@(9,'Vega Central Mapocho de Santiago','Metropolitana',45007,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',420,6000,6500,6202,'$/bandeja 2 kilos','Provincia de Curicó',3101,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',45006,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',400,6000,6500,6250,'$/bandeja 2 kilos','Provincia de Linares',3125,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44959,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',350,8000,8000,8000,'$/bandeja 2 kilos','Provincia de Curicó',4000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44196,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',550,6500,7000,6818,'$/bandeja 2 kilos','Provincia de Linares',3409,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44942,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',300,8000,8000,8000,'$/bandeja 2 kilos','Región de O''Higgins',4000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44294,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',480,7500,8000,7792,'$/bandeja 2 kilos','Provincia de Linares',3896,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44166,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Especial',50,8000,8000,8000,'$/bandeja 2 kilos','Provincia de Linares',4000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44166,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Especial',150,7200,7200,7200,'$/bandeja 2 kilos','Región de O''Higgins',3600,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44166,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',80,7000,7000,7000,'$/bandeja 2 kilos','Provincia de Linares',3500,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44637,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Especial',280,8000,8000,8000,'$/bandeja 2 kilos','Provincia de Linares',4000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44637,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',350,7000,7000,7000,'$/bandeja 2 kilos','Provincia de Linares',3500,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44922,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',380,8000,8000,8000,'$/bandeja 2 kilos','Provincia de Curicó',4000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44356,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',60,10000,10000,10000,'$/bandeja 2 kilos','Provincia de Curicó',5000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44222,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',200,6500,7000,6800,'$/bandeja 2 kilos','Provincia de Linares',3400,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44566,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Especial',410,8000,8000,8000,'$/bandeja 2 kilos','Provincia de Curicó',4000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44566,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',450,7000,7000,7000,'$/bandeja 2 kilos','Provincia de Curicó',3500,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44658,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',480,8000,8000,8000,'$/bandeja 2 kilos','Provincia de Linares',4000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44273,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',210,6000,6000,6000,'$/bandeja 2 kilos','Provincia de Linares',3000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44363,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',50,10000,10000,10000,'$/bandeja 2 kilos','Provincia de Curicó',5000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44974,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',300,7000,7000,7000,'$/bandeja 2 kilos','Provincia de Curicó',3500,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44974,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Segunda',250,6000,6000,6000,'$/bandeja 2 kilos','Provincia de Curicó',3000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44195,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',408,6509,7000,6774,'$/bandeja 2 kilos','Provincia de Linares',3387,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44194,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',190,5800,6000,5916,'$/bandeja 2 kilos','Provincia de Linares',2958,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44987,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',300,7000,7000,7000,'$/bandeja 2 kilos','Provincia de Curicó',3500,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44286,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',100,8000,8000,8000,'$/bandeja 2 kilos','Provincia de Linares',4000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44586,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Especial',250,8000,8000,8000,'$/bandeja 2 kilos','Provincia de Linares',4000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44586,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',200,7000,7000,7000,'$/bandeja 2 kilos','Provincia de Linares',3500,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44907,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',470,8000,8500,8266,'$/bandeja 2 kilos','Provincia de Curicó',4133,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44979,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',250,6000,6000,6000,'$/bandeja 2 kilos','Provincia de Curicó',3000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44979,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Segunda',150,5600,5600,5600,'$/bandeja 2 kilos','Provincia de Curicó',2800,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44901,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',400,7000,8000,7450,'$/bandeja 2 kilos','Región de O''Higgins',3725,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44650,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',580,7500,8000,7759,'$/bandeja 2 kilos','Provincia de Linares',3880,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44238,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',200,6000,6000,6000,'$/bandeja 2 kilos','Provincia de Curicó',3000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44238,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',150,6000,6000,6000,'$/bandeja 2 kilos','Provincia de Linares',3000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44165,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Especial',25,10000,10000,10000,'$/bandeja 2 kilos','Provincia de Linares',5000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44204,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',120,6500,6500,6500,'$/bandeja 2 kilos','Provincia de Linares',3250,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44204,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Segunda',150,7000,7000,7000,'$/bandeja 2 kilos','Provincia de Linares',3500,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44229,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',150,6000,6000,6000,'$/bandeja 2 kilos','Provincia de Curicó',3000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44663,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',450,8000,8000,8000,'$/bandeja 2 kilos','Provincia de Linares',4000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44938,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',290,8000,8000,8000,'$/bandeja 2 kilos','Región de O''Higgins',4000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44292,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',120,8000,8000,8000,'$/bandeja 2 kilos','Provincia de Linares',4000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44910,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',450,7000,7000,7000,'$/bandeja 2 kilos','Provincia de Curicó',3500,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44957,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',400,7000,7000,7000,'$/bandeja 2 kilos','Provincia de Curicó',3500,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44202,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',310,6500,7000,6677,'$/bandeja 2 kilos','Provincia de Linares',3338,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44651,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',380,8000,8000,8000,'$/bandeja 2 kilos','Provincia de Linares',4000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44673,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',250,8000,8000,8000,'$/bandeja 2 kilos','Provincia de Linares',4000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44568,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',680,7500,8000,7757,'$/bandeja 2 kilos','Provincia de Linares',3878,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44636,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',450,8000,8000,8000,'$/bandeja 2 kilos','Provincia de Linares',4000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44306,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',200,7000,7000,7000,'$/bandeja 2 kilos','Provincia de Curicó',3500,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',45008,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',570,6000,6500,6307,'$/bandeja 2 kilos','Provincia de Curicó',3154,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44960,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',350,8000,8000,8000,'$/bandeja 2 kilos','Región de O''Higgins',4000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44364,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',75,10000,10000,10000,'$/bandeja 2 kilos','Provincia de Curicó',5000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44215,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',280,6600,7000,6829,'$/bandeja 2 kilos','Provincia de Linares',3414,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44215,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Segunda',560,5600,6000,5814,'$/bandeja 2 kilos','Provincia de Linares',2907,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44209,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',370,5800,6000,5935,'$/bandeja 2 kilos','Provincia de Linares',2968,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44972,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',350,5600,5600,5600,'$/bandeja 2 kilos','Región de O''Higgins',2800,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44671,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',450,8000,8000,8000,'$/bandeja 2 kilos','Provincia de Linares',4000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44988,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',290,7000,7000,7000,'$/bandeja 2 kilos','Región Metropolitana',3500,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44236,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',450,6000,6000,6000,'$/bandeja 2 kilos','Provincia de Curicó',3000,2),
@(9,'Vega Central Mapocho de Santiago','Metropolitana',44210,13,'Fruta',100101,'Berries',100101004,'Frambuesa','Sin especificar','Primera',400,5800,6000,5910,'$/bandeja 2 kilos','Provincia de Linares',2955,2)
)

$nrows = $data.Count
$ncols = 20
$arr = New-Object 'object[,]' $nrows,$ncols
for ($i = 0; $i -lt $nrows; $i++) {
    for ($j = 0; $j -lt $ncols; $j++) {
        $arr[$i,$j] = $data[$i][$j]
    }
}

$ws.Range("A97:T156").Value = $arr

# Restore the date-time number format on column D for rows 97-156
# (matches existing formatting used for the other date cells in column D)
$ws.Range("D97:D156").NumberFormat = "YYYY-MM-DD HH:MM:SS"
